# Regenerate this week's report: bump the "Report Generated On" timestamp,
# populate the billed-amount / line-total figures now that the work request
# has real data (was a zero/#INVALID VALUE placeholder), and clear the
# Scope ID's "#INVALID VALUE" placeholder text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D5 - "Report Generated On:" timestamp refreshed for this regeneration run
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:01 AM"

# C8 - Total Billed Amount (Report Summary)
$ws.Range("C8").Value = 117

# G10 - Scope ID # placeholder had no real value; clear the "#INVALID VALUE" text
$ws.Range("G10").Value = ""

# H16 - Pricing for the single detail line item (Point 13 / CON-10-AAA-1-B-REEL)
$ws.Range("H16").Value = 117

# H17 - TOTAL row, mirrors the line item total above
$ws.Range("H17").Value = 117
